$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.300.58"
$ws.Range("E2").Value = "  -0.06%  "
$ws.Range("D3").Value = "'2.279.35"
$ws.Range("E3").Value = "  +1.02%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("D5").Value = "'503.65"
$ws.Range("E5").Value = "  +1.59%  "
$ws.Range("D6").Value = "'129.25"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.528"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'0.0954"
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "'0.333"
$ws.Range("E11").Value = "  +3.44%  "
$ws.Range("D12").Value = "'4.72"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").Value = "'2.687.08"
$ws.Range("E13").Value = "  +1.03%  "
$ws.Range("D14").Value = "'22.83"
$ws.Range("E14").Value = "  +5.56%  "
$ws.Range("D15").Value = "'54.255.67"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("E16").Value = "  +0.22%  "
$ws.Range("D17").Value = "'2.299.76"
$ws.Range("E17").Value = "  +0.82%  "
$ws.Range("D18").Value = "'10.25"
$ws.Range("E18").Value = "  +3.18%  "
$ws.Range("D19").Value = "'4.12"
$ws.Range("E19").Value = "  +1.93%  "
$ws.Range("D20").Value = "'303.91"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("D21").Value = "'6.42"
$ws.Range("E21").Value = "  +1.21%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("D23").Value = "'61.84"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("E26").Value = "  +2.97%  "
$ws.Range("D27").Value = "'174.66"
$ws.Range("E27").Value = "  +5.93%  "
$ws.Range("D28").Value = "'1.62"
$ws.Range("E28").Value = "  +1.58%  "
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("D30").Value = "'0.0₃0688"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  +0.46%  "
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'17.78"
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("B34").Value = "SuiNetwork"
$ws.Range("C34").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D34").Value = "'0.954"
$ws.Range("E34").Value = "  +9.55%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "'0.997"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("E36").Value = "  +1.02%  "
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  +2.86%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  +0.69%  "
$ws.Range("E40").Value = "  +1.46%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'4.82"
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'125.07"
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("E43").Value = "  +3.42%  "
$ws.Range("D44").Value = "'0.0893"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "'0.547"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").Value = "'240.21"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("E48").Value = "  +0.80%  "
$ws.Range("D49").Value = "'10.77"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "'16.40"
$ws.Range("E50").Value = "  +0.56%  "
$ws.Range("E51").Value = "  +0.17%  "
